$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; F = 1; G = 6.189590430959694 }
    3 = @{ B = 0.6606524410359556; C = 0.04071648406533734; D = 0.1494219747398047; E = 0.4942365360607697; F = 1; G = 1.345027435901867 }
    4 = @{ B = 1.455362044514542;  C = 1.655778082260271;  D = 0.7527432677738641; E = 0.4942365360607697; F = 1; G = 4.358119930609447 }
    5 = @{ B = 0.6606524410359556; C = 10.34677158129881;  D = 3.537761648806719;  E = 10.19245300693656;  F = 1; G = 24.73763867807805 }
    6 = @{ B = 1.455362044514542;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697; F = 0; G = 3.754798637575387 }
    7 = @{ B = 0.1190320826869504; C = 0.306821227259698;  D = 0.7527432677738641; E = 0.4942365360607697; F = 1; G = 1.672833113781282 }
    8 = @{ B = 1.455362044514542;  C = 1.655778082260271;  D = 22.3905356188092;   E = 0.4942365360607697; F = 1; G = 25.99591228164478 }
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
